$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# "일지 (임시)" (work log) sheet: append 6 new journal rows (12-17)
# ------------------------------------------------------------------
$wsLog = $wb.Worksheets.Item("일지 (임시)")

# Row 12 - 2018/08/26, 김대홍, ovenapp 프로토타이핑 작성
$wsLog.Cells.Item(5, 1).Copy($wsLog.Cells.Item(12, 1))
$wsLog.Cells.Item(12, 1).Value2 = 43338
$wsLog.Cells.Item(12, 2).Value = "김대홍"
$wsLog.Cells.Item(12, 3).Value = "ovenapp 프로토타이핑 작성"

# Row 13 - 2018/08/27, 김대홍, 현황 분석
$wsLog.Cells.Item(5, 1).Copy($wsLog.Cells.Item(13, 1))
$wsLog.Cells.Item(13, 1).Value2 = 43339
$wsLog.Cells.Item(13, 2).Value = "김대홍"
$wsLog.Cells.Item(13, 3).Value = "현황 분석"

# Row 14 - 2018/08/27, 김대홍, 현 컨텐츠 기준으로 구현할 메뉴 정리 및 회의
$wsLog.Cells.Item(5, 1).Copy($wsLog.Cells.Item(14, 1))
$wsLog.Cells.Item(14, 1).Value2 = 43339
$wsLog.Cells.Item(14, 2).Value = "김대홍"
$wsLog.Cells.Item(14, 3).Value = "현 컨텐츠 기준으로 구현할 메뉴 정리 및 회의"

# Row 15 - 2018/08/28, 김대홍, DB 환경 구축
$wsLog.Cells.Item(5, 1).Copy($wsLog.Cells.Item(15, 1))
$wsLog.Cells.Item(15, 1).Value2 = 43340
$wsLog.Cells.Item(15, 2).Value = "김대홍"
$wsLog.Cells.Item(15, 3).Value = "DB 환경 구축"

# Row 16 - 2018/08/29, 김대홍, joinform, loginform 틀 구현
$wsLog.Cells.Item(5, 1).Copy($wsLog.Cells.Item(16, 1))
$wsLog.Cells.Item(16, 1).Value2 = 43341
$wsLog.Cells.Item(16, 2).Value = "김대홍"
$wsLog.Cells.Item(16, 3).Value = "joinform, loginform 틀 구현"

# Row 17 - 2018/08/29, 김대홍, TB_USER_INFO 테이블 및 컬럼 생성
$wsLog.Cells.Item(5, 1).Copy($wsLog.Cells.Item(17, 1))
$wsLog.Cells.Item(17, 1).Value2 = 43341
$wsLog.Cells.Item(17, 2).Value = "김대홍"
$wsLog.Cells.Item(17, 3).Value = "TB_USER_INFO 테이블 및 컬럼 생성"

# ------------------------------------------------------------------
# Restore the cursor/selection on a few sheets that were browsed
# while the author was working, then leave "일지 (임시)" active on
# its new last row of data (as the final tab the author left open).
# ------------------------------------------------------------------
$wsClassify = $wb.Worksheets.Item("분류 설명")
$wsClassify.Activate() | Out-Null
$wsClassify.Range("F16").Select() | Out-Null

$wsDb = $wb.Worksheets.Item("DB 구조")
$wsDb.Activate() | Out-Null
$wsDb.Range("B9").Select() | Out-Null

$wsTable = $wb.Worksheets.Item("TABLE 구조")
$wsTable.Activate() | Out-Null
$wsTable.Range("A5").Select() | Out-Null

$wsLog.Activate() | Out-Null
$wsLog.Range("C19").Select() | Out-Null
